# Update the embedded "version" string that is stamped throughout the
# workbook from the old build ("mines - version 1.0.0 (Feb 3 2026) (built on
# February 03 2026 10.14.00 EST)") to the new release build
# ("Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on
# February 03 2026 17.29.55 EST)").

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet -------------------------------------------------------

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Ensham Coal Mine, Australia, M0038, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# --- "Boundaries and methane sources" sheet ------------------------------

for ($row = 2; $row -le 9; $row++) {
    $wsData.Range("S$row").Value = $newVersion
}
